$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.302.86"
$ws.Range("E2").Value = "  +3.82%  "
$ws.Range("D3").Value = "3.124.40"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("D5").Value = "'219.89"
$ws.Range("E5").Value = "  +5.00%  "
$ws.Range("D6").Value = "'625.61"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").Value = "'0.383"
$ws.Range("E7").Value = "  +3.32%  "
$ws.Range("D8").Value = "'0.969"
$ws.Range("E8").Value = "  +19.34%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "3.123.87"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("E11").Value = "  +17.96%  "
$ws.Range("E12").Value = "  +5.99%  "
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  +5.48%  "
$ws.Range("D14").Value = "'34.69"
$ws.Range("E14").Value = "  +9.03%  "
$ws.Range("D15").Value = "'5.40"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").Value = "91.196.15"
$ws.Range("E16").Value = "  +3.91%  "
$ws.Range("D17").Value = "3.704.20"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "3.124.78"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").Value = "'3.76"
$ws.Range("E19").Value = "  +17.16%  "
$ws.Range("D20").Value = "'0.0000219"
$ws.Range("E20").Value = "  +9.19%  "
$ws.Range("D21").Value = "'14.13"
$ws.Range("E21").Value = "  +7.05%  "
$ws.Range("D22").Value = "'437.47"
$ws.Range("E22").Value = "  +4.24%  "
$ws.Range("D23").Value = "'8.76"
$ws.Range("E23").Value = "  +7.89%  "
$ws.Range("E24").Value = "  +6.32%  "
$ws.Range("D25").Value = "'6.13"
$ws.Range("E25").Value = "  +12.03%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'86.49"
$ws.Range("E26").Value = "  +5.65%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "'12.22"
$ws.Range("E27").Value = "  +3.75%  "
$ws.Range("D28").Value = "3.295.09"
$ws.Range("E28").Value = "  +2.26%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'0.169"
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("E31").Value = "  +11.98%  "
$ws.Range("D32").Value = "'528.23"
$ws.Range("E32").Value = "  +4.77%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "'0.899"
$ws.Range("E33").Value = "  -16.08%  "
$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Value = "'3.77"
$ws.Range("E34").Value = "  +6.56%  "
$ws.Range("D35").Value = "'7.13"
$ws.Range("E35").Value = "  +7.73%  "
$ws.Range("E36").Value = "  +9.17%  "
$ws.Range("D37").Value = "'23.78"
$ws.Range("E37").Value = "  +7.09%  "
$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").Value = "'1.87"
$ws.Range("E38").Value = "  +3.96%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'1.28"
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("D40").Value = "'22.29"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +15.85%  "
$ws.Range("D43").Value = "'0.0803"
$ws.Range("E43").Value = "  +20.10%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("B45").Value = "PolygonEcosystemToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D45").Value = "'0.380"
$ws.Range("E45").Value = "  +5.88%  "
$ws.Range("D46").Value = "'1.92"
$ws.Range("E46").Value = "  +7.12%  "
$ws.Range("D47").Value = "'146.14"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").Value = "'44.20"
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("D49").Value = "'1.30"
$ws.Range("E49").Value = "  +10.61%  "
$ws.Range("D50").Value = "'167.35"
$ws.Range("E50").Value = "  +7.89%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "'4.19"
$ws.Range("E51").Value = "  +7.15%  "
